# Update "想去人数" (interest count) figures pulled at the latest gh-pages
# data refresh (build 456a3b4).
$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 3697
$wsExpo.Range("F6").Value = 38
$wsExpo.Range("F7").Value = 191

# 演出 (Performances) sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 123

# 全部类型 (All types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 123
$wsAll.Range("F8").Value = 3697
$wsAll.Range("F10").Value = 38
$wsAll.Range("F12").Value = 191
